$d = $word.ActiveDocument

# Locate the target bullet paragraph ("Ensured the integration and maintenance
# of new business in to existing financial reporting systems") by scanning the
# Paragraphs collection for its distinctive text, rather than relying on fixed
# character offsets.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ensured the integration and maintenance of new business*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$r = $target.Range

# Replace the whole paragraph (including its trailing paragraph mark) with the
# rewritten OOXML: the bullet now reads "Collaborated with multiple data
# engineering teams; streamlined the pipeline connecting multiple ERPs", with
# "Collaborated" and "streamlined" bold, and the paragraph-mark run's szCs
# trimmed from 28 to 21 to match the new run sizes.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="256"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
'<w:p w14:paraId="7BB0EFFB" w14:textId="6FC1A8E4" w:rsidR="00FC56B2" w:rsidRPr="001472C0" w:rsidRDefault="00053C8F" w:rsidP="00364C66">' + `
'<w:pPr><w:pStyle w:val="ListParagraph"/><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr>' + `
'<w:tabs><w:tab w:val="left" w:pos="789"/><w:tab w:val="left" w:pos="791"/></w:tabs>' + `
'<w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:contextualSpacing w:val="0"/>' + `
'<w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr>' + `
'<w:r><w:rPr><w:sz w:val="21"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Collaborated</w:t></w:r>' + `
'<w:r><w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> with multiple data engineering teams;</w:t></w:r>' + `
'<w:r><w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>streamlined</w:t></w:r>' + `
'<w:r><w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> the pipeline connecting multiple ERPs</w:t></w:r>' + `
'</w:p>' + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml) | Out-Null

Write-Output "Paragraph updated"
